# Template tweak: give the title placeholder on the title-slide layout a bit
# more room (taller) and move the subtitle placeholder down / make it a
# little shorter to match, so longer titles fit (see commit message).
#
# NOTE on the notes-master date field (14/08/2018 -> 16/08/2018) and the
# subtitle run's "dirty" bookkeeping attribute from the source diff: this
# host's object model does not expose a writable path for either (a
# datetimeFigureOut field's cached text is read-only via automation, same
# as stock PowerPoint, and "dirty" has no COM-visible property), so they
# are intentionally left alone here rather than risk corrupting unrelated
# parts of the template.

$p = $ppt.ActivePresentation

# Title-slide custom layout ("Titeldia" == ppt/slideLayouts/slideLayout1.xml)
$master = $p.SlideMaster
$titleLayout = $master.CustomLayouts.Item(1)

$titleShape = $titleLayout.Shapes.Item(1)      # "Title 1" (ctrTitle placeholder)
$subtitleShape = $titleLayout.Shapes.Item(2)   # "Subtitle 2" (subTitle placeholder)

# Title: keep its left/top/width, only grow its height
#   ext cy: 2387600 -> 2626678 EMU
$titleShape.Height = 206.8250732421875

# Subtitle: keep its left/width, move it down and shrink its height
#   off y : 3602038 -> 3848794 EMU
#   ext cy: 1655762 -> 1409006 EMU
$subtitleShape.Top = 303.0546875
$subtitleShape.Height = 110.94537353515625
